# mirroring.pptx "Initial commit for kfserving abtesting" edit, replayed
# against the single-slide reduced deck.
#
# EMU -> point conversion used throughout: 1 pt = 12700 EMU (PowerPoint's
# COM surface works in points for Left/Top/Width/Height).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# 1) Resize "Rectangle 10" (the big white background panel): it moves down
#    and shrinks a bit (the header band above it got taller).
$rect10 = $s.Shapes.Item("Rectangle 10")
$rect10.Left   = 0 / 12700
$rect10.Top    = 1389527 / 12700
$rect10.Width  = 12191999 / 12700
$rect10.Height = 4518745 / 12700

# 2) Delete the small empty "TextBox 9" placeholder textbox near the top.
$s.Shapes.Item("TextBox 9").Delete()

# 3) Delete "TextBox 4" ("Istio Virtual Service" label).
$s.Shapes.Item("TextBox 4").Delete()

# 4) Delete "TextBox 95" ("winner" label).
$s.Shapes.Item("TextBox 95").Delete()

# 5) Delete the standalone "Ribbon" icon picture ("Graphic 74").
$s.Shapes.Item("Graphic 74").Delete()

# 6) Delete "Rounded Rectangle 106" ("v2.0" badge).
$s.Shapes.Item("Rounded Rectangle 106").Delete()

# 7) Delete "TextBox 124" ("Objectives" label).
$s.Shapes.Item("TextBox 124").Delete()

# 8) Delete the "Checkbox Checked" icon picture ("Graphic 125").
$s.Shapes.Item("Graphic 125").Delete()

# 9) & 10) Remove the footer and slide-number placeholders. Turning off
#     visibility (rather than Shape.Delete, which would just respawn a
#     fresh placeholder instance from the layout) removes the shapes
#     outright, matching the target deck which has neither shape left.
$hf = $s.HeadersFooters
$hf.Footer.Visible = $false
$hf.SlideNumber.Visible = $false
